$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newParticipantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.instrument_model in ['Illumina HiSeq 2000']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@
$newSamplesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,apoc.coll.sort(COLLECT(DISTINCT samp.sample_tumor_status)) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@
$newFilesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2000']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,apoc.coll.sort(COLLECT(DISTINCT samp.sample_tumor_status)) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@

# Update B3 (SamplesTab query) to the new samples query (with apoc.coll.sort)
$ws.Range("B3").Value2 = $newSamplesQuery

# Update B4 (FilesTab query) to the new files query (with apoc.coll.sort)
$ws.Range("B4").Value2 = $newFilesQuery

# Update B2 (ParticipantsTab query) to the new participants query
$ws.Range("B2").Value2 = $newParticipantsQuery

# Adjust row heights to match the new content (auto-fit would normally occur in Excel on edit)
$ws.Rows(2).RowHeight = 354.75
$ws.Rows(3).RowHeight = 282.75

# Update the sheet view selection to C2 (also clears the previous scrolled topLeftCell)
$ws.Range("C2").Select()
